# The edit shuffles the "Fecha" (and corresponding Volumen/Precio/Origen/etc.)
# data across rows 3-16 of the sheet (row 9 is left untouched). Effectively,
# each target row ends up with the full row-contents that used to live in a
# different source row. Capture all the old row values first (so we don't
# clobber data we still need to read), then write them back out in the new
# order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row (values to copy from the original sheet)
$rowMap = @{
    3  = 15
    4  = 5
    5  = 4
    6  = 11
    7  = 14
    8  = 6
    10 = 12
    11 = 8
    12 = 3
    13 = 16
    14 = 10
    15 = 7
    16 = 13
}

$firstRow = 3
$lastRow = 16
$firstCol = 1   # A
$lastCol = 18   # R

# Snapshot the original values for every row we might need as a source.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowValues = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowValues[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowValues
}

# Write back each destination row using the snapshot of its source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcValues = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcValues[$c]
    }
}
